$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column D (current D:K quarterly data shifts to E:L)
$ws.Range("D1").EntireColumn.Insert()

# Copy number formats/styles from the (now shifted) column E into the new column D
# so the new quarter's column matches the formatting of its neighboring data (dates vs numbers).
# Done per contiguous data block, skipping the blank separator rows (36, 78) and the
# section-header rows (5, 6, 37, 79) that only carry a label in column A/B.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column D with the latest quarter's figures
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 400
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 600
$ws.Range("D18").Value = -200
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = -300
$ws.Range("D22").Value = "NA"
$ws.Range("D23").Value = -300
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -300
$ws.Range("D27").Value = -300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = -300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -300
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 400
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 100
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 3700
$ws.Range("D48").Value = 0
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 4200
$ws.Range("D57").Value = 100
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = "NA"
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 100
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -32400
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 4200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = -300
$ws.Range("D83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -200
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 200
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 0

Write-Host "Inserted newest quarter column and populated values"
